# Generate Report for Archive
#
# 1. Update status text "Ready for handoff" -> "In Translation" wherever it
#    appears (Overview sheet E2/F2, and the Status column (C2) on the
#    per-locale "zh-cn" and "de-de" sheets).
# 2. Narrow the "Status" column(s) from ~17.22 chars to ~13.41 chars:
#    - Overview sheet: columns E and F (zh-cn / de-de status columns)
#    - zh-cn sheet: column C (Status)
#    - de-de sheet: column C (Status)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text update: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width updates (narrower Status columns) ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
